# Trade #71 closed at 2026-02-17 12:54:30 - unknown UNKNOWN +0.000%
#
# Updates summary / strategy-status counters and appends the new trade
# row (#71, 1-indexed Trade # = 71, sheet row 72) to both the "All Trades"
# and "MarketMaking" worksheets.

$wb = $excel.ActiveWorkbook

# --- Summary sheet -----------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B6").Value = 71
$wsSummary.Range("B9").Value = 45.07

# --- Strategy Status sheet ---------------------------------------------
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("D4").Value = 71
$wsStatus.Range("G4").Value = 45.07

# --- Append new trade row to both trade-log sheets ----------------------
$newRow = 72

$tradeNum       = 71
$tradeDate      = "2026-02-17"
$tradeTime      = "12:54:24"
$tradeStrategy  = "MarketMaking"
$tradeSide      = "DOWN"
$entryPrice     = 0.935106
$exitPrice      = 0.9399999999999999
$status         = "CLOSED"
$pnlPct         = 0.5234
$pnlDollar      = 0
$capitalAfter   = 100.21
$entrySlippage  = 0
$exitSlippage   = 0
$confidence     = 0.6
$entryReason    = "Normal spread capture: 19600 bps"
$exitReason     = "early_exit"
$durationMin    = 0.13

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Cells.Item($newRow, 1).Value  = $tradeNum
    $ws.Cells.Item($newRow, 2).Value  = $tradeDate
    $ws.Cells.Item($newRow, 3).Value  = $tradeTime
    $ws.Cells.Item($newRow, 4).Value  = $tradeStrategy
    $ws.Cells.Item($newRow, 5).Value  = $tradeSide
    $ws.Cells.Item($newRow, 6).Value  = $entryPrice
    $ws.Cells.Item($newRow, 7).Value  = $exitPrice
    $ws.Cells.Item($newRow, 8).Value  = $status
    $ws.Cells.Item($newRow, 9).Value  = $pnlPct
    $ws.Cells.Item($newRow, 10).Value = $pnlDollar
    $ws.Cells.Item($newRow, 11).Value = $capitalAfter
    $ws.Cells.Item($newRow, 12).Value = $entrySlippage
    $ws.Cells.Item($newRow, 13).Value = $exitSlippage
    $ws.Cells.Item($newRow, 14).Value = $confidence
    $ws.Cells.Item($newRow, 15).Value = $entryReason
    $ws.Cells.Item($newRow, 16).Value = $exitReason
    $ws.Cells.Item($newRow, 17).Value = $durationMin
}
